$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at 462, shifting old rows 462-472 down to 468-478
$ws.Rows("462:467").Insert()

# Populate the new rows (462-467) with the new date block (2022-02-03 / serial 44595)
# Row 462
$ws.Cells.Item(462, 1).Value = 2
$ws.Cells.Item(462, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(462, 3).Value = "Coquimbo"
$ws.Cells.Item(462, 4).Value = 44595
$ws.Cells.Item(462, 5).Value = 4
$ws.Cells.Item(462, 6).Value = 100112020
$ws.Cells.Item(462, 7).Value = "Tomate"
$ws.Cells.Item(462, 8).Value = "Larga vida"
$ws.Cells.Item(462, 9).Value = "Primera"
$ws.Cells.Item(462, 10).Value = 1400
$ws.Cells.Item(462, 11).Value = 10000
$ws.Cells.Item(462, 12).Value = 11000
$ws.Cells.Item(462, 13).Value = 10500
$ws.Cells.Item(462, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(462, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(462, 16).Value = 583
$ws.Cells.Item(462, 17).Value = 18
$ws.Cells.Item(462, 18).Value = "Hortaliza"

# Row 463
$ws.Cells.Item(463, 1).Value = 2
$ws.Cells.Item(463, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(463, 3).Value = "Coquimbo"
$ws.Cells.Item(463, 4).Value = 44595
$ws.Cells.Item(463, 5).Value = 4
$ws.Cells.Item(463, 6).Value = 100112020
$ws.Cells.Item(463, 7).Value = "Tomate"
$ws.Cells.Item(463, 8).Value = "Larga vida"
$ws.Cells.Item(463, 9).Value = "Segunda"
$ws.Cells.Item(463, 10).Value = 800
$ws.Cells.Item(463, 11).Value = 8000
$ws.Cells.Item(463, 12).Value = 9000
$ws.Cells.Item(463, 13).Value = 8500
$ws.Cells.Item(463, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(463, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(463, 16).Value = 472
$ws.Cells.Item(463, 17).Value = 18
$ws.Cells.Item(463, 18).Value = "Hortaliza"

# Row 464
$ws.Cells.Item(464, 1).Value = 2
$ws.Cells.Item(464, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(464, 3).Value = "Coquimbo"
$ws.Cells.Item(464, 4).Value = 44595
$ws.Cells.Item(464, 5).Value = 4
$ws.Cells.Item(464, 6).Value = 100112020
$ws.Cells.Item(464, 7).Value = "Tomate"
$ws.Cells.Item(464, 8).Value = "Larga vida"
$ws.Cells.Item(464, 9).Value = "Tercera"
$ws.Cells.Item(464, 10).Value = 400
$ws.Cells.Item(464, 11).Value = 6000
$ws.Cells.Item(464, 12).Value = 7000
$ws.Cells.Item(464, 13).Value = 6500
$ws.Cells.Item(464, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(464, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(464, 16).Value = 361
$ws.Cells.Item(464, 17).Value = 18
$ws.Cells.Item(464, 18).Value = "Hortaliza"

# Row 465
$ws.Cells.Item(465, 1).Value = 2
$ws.Cells.Item(465, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(465, 3).Value = "Coquimbo"
$ws.Cells.Item(465, 4).Value = 44595
$ws.Cells.Item(465, 5).Value = 4
$ws.Cells.Item(465, 6).Value = 100112020
$ws.Cells.Item(465, 7).Value = "Tomate"
$ws.Cells.Item(465, 8).Value = "Semiduro"
$ws.Cells.Item(465, 9).Value = "Primera"
$ws.Cells.Item(465, 10).Value = 1800
$ws.Cells.Item(465, 11).Value = 6500
$ws.Cells.Item(465, 12).Value = 7000
$ws.Cells.Item(465, 13).Value = 6750
$ws.Cells.Item(465, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(465, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(465, 16).Value = 375
$ws.Cells.Item(465, 17).Value = 18
$ws.Cells.Item(465, 18).Value = "Hortaliza"

# Row 466
$ws.Cells.Item(466, 1).Value = 2
$ws.Cells.Item(466, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(466, 3).Value = "Coquimbo"
$ws.Cells.Item(466, 4).Value = 44595
$ws.Cells.Item(466, 5).Value = 4
$ws.Cells.Item(466, 6).Value = 100112020
$ws.Cells.Item(466, 7).Value = "Tomate"
$ws.Cells.Item(466, 8).Value = "Semiduro"
$ws.Cells.Item(466, 9).Value = "Segunda"
$ws.Cells.Item(466, 10).Value = 1000
$ws.Cells.Item(466, 11).Value = 4500
$ws.Cells.Item(466, 12).Value = 5000
$ws.Cells.Item(466, 13).Value = 4750
$ws.Cells.Item(466, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(466, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(466, 16).Value = 264
$ws.Cells.Item(466, 17).Value = 18
$ws.Cells.Item(466, 18).Value = "Hortaliza"

# Row 467
$ws.Cells.Item(467, 1).Value = 2
$ws.Cells.Item(467, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(467, 3).Value = "Coquimbo"
$ws.Cells.Item(467, 4).Value = 44595
$ws.Cells.Item(467, 5).Value = 4
$ws.Cells.Item(467, 6).Value = 100112020
$ws.Cells.Item(467, 7).Value = "Tomate"
$ws.Cells.Item(467, 8).Value = "Semiduro"
$ws.Cells.Item(467, 9).Value = "Tercera"
$ws.Cells.Item(467, 10).Value = 400
$ws.Cells.Item(467, 11).Value = 2500
$ws.Cells.Item(467, 12).Value = 3000
$ws.Cells.Item(467, 13).Value = 2750
$ws.Cells.Item(467, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(467, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(467, 16).Value = 153
$ws.Cells.Item(467, 17).Value = 18
$ws.Cells.Item(467, 18).Value = "Hortaliza"
